$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Replace the Q&A corpus in rows 2-3 (existing rows) ---
$ws.Range("A2").Value = "What is the University Scholars Programme?"
$ws.Range("B2").Value = "The University Scholars Programme (USP) is an undergraduate academic programme established in 2001 in the National University of Singapore (NUS)."

$ws.Range("A3").Value = "What is USP?"
$ws.Range("B3").Value = "The University Scholars Programme (USP) is an undergraduate academic programme established in 2001 in the National University of Singapore (NUS)."

# --- Append the rest of the new FAQ corpus ---
$ws.Range("A4").Value = "What are USP's values?"
$ws.Range("B4").Value = "A USP student values intense academic inquiry, research, exploration and engagement."

$ws.Range("A5").Value = "What does it mean to be curious?"
$ws.Range("B5").Value = "Curious means adventurous in thought, pursuing a broad range of intellectual interests and ready to make connections across different domains."

$ws.Range("A6").Value = "What does it mean to be critical?"
$ws.Range("B6").Value = "Critical means intellectually rigorous, deeply reflective and having a humility born of awareness of the limitations of our own knowledge."

$ws.Range("A7").Value = "What does it mean to be courageous?"
$ws.Range("B7").Value = "Courageous means willing to consider differing points of view, unafraid to face challenges and to act upon ideas."

$ws.Range("A8").Value = "What does it mean to be engaged?"
$ws.Range("B8").Value = "Engaged means extending the frontiers of knowledge in energetic and creative ways, prepared to navigate and help shape a complex world in a responsible way."

$ws.Range("A9").Value = "What faculties are there in USP?"
$ws.Range("B9").Value = "Students admitted to USP are concurrently enrolled in 1 of 7 NUS faculties or schools: Faculty of Arts and Social Sciences, Faculty of Engineering, Faculty of Science, Faculty of Law, NUS Business School, School of Computing and School of Design and Environment."

$ws.Range("A10").Value = "What are the academic requirements of USP?"
$ws.Range("B10").Value = "You should, in general, complete your degree with honours within four years. As a student in USP, you will have to read and pass the following: 3 Foundation Tier modules, 8 Inquiry Tier modules and 1 Reflection Tier module. Please refer to the respective cohorts requirements for details: http://www.usp.nus.edu.sg/curriculum/academic-requirements/"

$ws.Range("A11").Value = "How many students are admitted to USP?"
$ws.Range("B11").Value = "Each year, about 200 incoming NUS undergraduates are admitted to USP."

$ws.Range("A12").Value = "How can I contact USP?"
$ws.Range("B12").Value = "USP Contacts: 18 College Avenue East, Singapore 138593, +65 6516 4425, General Enquiries: usphelp@nus.edu.sg"

# --- Formatting: the long "faculties" answer gets vertical-center alignment ---
$ws.Range("B9").Font.Name = "Calibri"
$ws.Range("B9").VerticalAlignment = -4108

# --- Column widths (auto-fit to the new, much wider text) ---
$ws.Columns.Item(1).ColumnWidth = 36.6
$ws.Columns.Item(2).ColumnWidth = 124.42

# --- Clear the logbook: move the selection/cursor past the new data ---
$ws.Range("A13").Select()

# --- Restore window placement/sizing recorded in the session ---
$excel.ActiveWindow.Left = -11220
$excel.ActiveWindow.Top = 5952
$excel.ActiveWindow.Width = 17280
$excel.ActiveWindow.Height = 8964
